$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1) Merge the split runs "Tamayo, " / "Aaron John" / " N." into a single run
#    with the text "Tamayo, Aaron John N." (Find/Replace across the whole
#    match collapses the runs it touches into one run in this engine).
# --------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Tamayo, Aaron John N.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tamayo, Aaron John N.", 2) | Out-Null

# --------------------------------------------------------------------------
# 2) Merge the split runs "Valencia, " / "Anton " / "Philip" into a single
#    run with the text "Valencia, Anton Philip".
# --------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Valencia, Anton Philip", $true, $false, $false, $false, $false,
    $true, 1, $false, "Valencia, Anton Philip", 2) | Out-Null

# --------------------------------------------------------------------------
# 3) Locate the "[150-600 word summary ...]" paragraph, add
#    <w:lang w:val="en-PH"/> to its paragraph-mark run properties (pPr/rPr),
#    and insert a brand new paragraph right after it containing
#    "ABARRA: Hello World!" (also tagged en-PH), matching the diff exactly.
# --------------------------------------------------------------------------
$summaryPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match [regex]::Escape("[150-600 word summary of the report that provides a high-level overview of the project]")) {
        $summaryPara = $p
        break
    }
}

if ($summaryPara -ne $null) {
    $summaryXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="00000018" w14:textId="77777777" w:rsidR="00961C1E" w:rsidRDefault="00000000"><w:pPr><w:spacing w:after="200"/><w:ind w:right="-41"/><w:rPr><w:rFonts w:ascii="Archivo" w:hAnsi="Archivo" w:cs="Archivo"/><w:lang w:val="en-PH"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Archivo" w:hAnsi="Archivo" w:cs="Archivo"/></w:rPr><w:t>[150-600 word summary of the report that provides a high-level overview of the project]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $summaryPara.Range.InsertXML($summaryXml)

    # Re-fetch the paragraph (range positions are unaffected by the XML swap
    # since it kept the same text), then add the new paragraph after it.
    $summaryPara.Range.InsertParagraphAfter()

    # The newly inserted (still empty) paragraph now directly follows.
    $newPara = $summaryPara.Next()

    $newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="200"/><w:ind w:right="-41"/><w:rPr><w:rFonts w:ascii="Archivo" w:hAnsi="Archivo" w:cs="Archivo"/><w:lang w:val="en-PH"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Archivo" w:hAnsi="Archivo" w:cs="Archivo"/><w:lang w:val="en-PH"/></w:rPr><w:t>ABARRA: Hello World!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $newPara.Range.InsertXML($newParaXml)
}
